# Refresh the cryptos table with the latest scrape (GitHub Actions run).
# Updates Price (D) and Volume(1h) (E) for most rows, and additionally
# swaps the Chainlink/Polkadot rows (18/19) to reflect the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.110.03"
$ws.Range("E2").Value = "  +1.98%  "
$ws.Range("D3").Value = "3.933.34"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.04%  "
$ws.Range("D7").Value = "3.933.08"
$ws.Range("E7").Value = "  +2.60%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("E12").Value = "  +2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000257"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.79%  "
$ws.Range("D15").Value = "4.594.38"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "3.888.49"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "70.124.97"
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.59%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "494.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.750"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.61%  "
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("E26").Value = "  +2.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.04%  "
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.85%  "
$ws.Range("D32").Value = "4.085.99"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "3.897.32"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("E38").Value = "  +1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.141"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.26%  "
$ws.Range("E40").Value = "  +12.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.331"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "439.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000278"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +24.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0372"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
